$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "QC" column (D), shifting QC/QT from D/E to E/F.
# Excel's Insert() mirrors the style of the (new) neighbouring column C onto the
# inserted cells, which already matches the header/body styles used by the target.
$ws.Columns("D").Insert()

# Header + sample note text for the newly inserted "Note" column.
$ws.Range("D1").Value = "Note"
$ws.Range("D2").Value = "The cat sat in front of the bird cage in an agony of frustration at being so near and yet so far.`n猫无可奈何地坐在鸟笼前，眼看着鸟儿近在咫尺，可怎么也够不着。"

# Match the authored column width (closest value reachable through the
# character-width COM API).
$ws.Columns("D").ColumnWidth = 60

# The AutoFilter range does not auto-grow when a column is inserted inside it,
# so turn it off and re-apply it across the new A1:F1 extent.
$ws.AutoFilterMode = $false
$reapplied = $ws.Range("A1:F1").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "animal!_FilterDatabase") {
        $n.RefersTo = "=animal!`$A`$1:`$F`$1"
    }
}
